$wb = $excel.ActiveWorkbook

function Set-Row($ws, $row, $vals) {
    for ($i = 0; $i -lt $vals.Count; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $vals[$i]
    }
}

# --- Sheet: Home win ---
$ws = $wb.Worksheets.Item("Home win")
Set-Row $ws 2 @("28-01-2025 19:45", "ENGLAND", "LEAGUE TWO", "Swindon Town - Tranmere", 70, 1.95)
Set-Row $ws 3 @("28-01-2025 19:45", "ENGLAND", "NATIONAL LEAGUE - NORTH", "Alfreton Town - Leamington", 80, 2.2)
Set-Row $ws 4 @("28-01-2025 08:00", "INDONESIA", "LIGA 2", "Persikas - Persipa Pati", 73.3, 3.9)
Set-Row $ws 5 @("28-01-2025 23:30", "WORLD", "SUDAMERICANO U20", "Colombia U20 - Ecuador U20", 73.3, 2.2)
Set-Row $ws 6 @("29-01-2025 22:00", "BRAZIL", "GAÚCHO - 1", "Ypiranga-RS - São Luiz", 86.7, 2.25)
Set-Row $ws 7 @("29-01-2025 19:45", "ENGLAND", "NATIONAL LEAGUE CUP", "Braintree - Tottenham Hotspur U21", 70, 2.25)

# --- Sheet: Away Win --- (unchanged per diff; left untouched)

# --- Sheet: Draw ---
$ws = $wb.Worksheets.Item("Draw")
Set-Row $ws 2 @("28-01-2025 20:00", "ENGLAND", "LEAGUE ONE", "Huddersfield - Birmingham", 60, 3.4)
Set-Row $ws 3 @("28-01-2025 15:00", "MEXICO", "U23 LEAGUE", "Necaxa U23 - Cruz Azul U23", 60, 3.3)
Set-Row $ws 4 @("29-01-2025 22:30", "BRAZIL", "PAULISTA - A1", "Noroeste - Inter De Limeira", 80, 3)

# --- Sheet: Btts ---
$ws = $wb.Worksheets.Item("Btts")
Set-Row $ws 2 @("28-01-2025 19:45", "ENGLAND", "LEAGUE ONE", "Rotherham - Cambridge United", 80, 2)
Set-Row $ws 3 @("29-01-2025 20:00", "WORLD", "UEFA CHAMPIONS LEAGUE", "Aston Villa - Celtic", 80, 1.7)
Set-Row $ws 4 @("29-01-2025 20:00", "WORLD", "UEFA CHAMPIONS LEAGUE", "Bayern München - Slovan Bratislava", 80, 3)
Set-Row $ws 5 @("29-01-2025 20:00", "WORLD", "UEFA CHAMPIONS LEAGUE", "Inter - Monaco", 76.7, 1.7)
Set-Row $ws 6 @("29-01-2025 20:00", "WORLD", "UEFA CHAMPIONS LEAGUE", "Lille - Feyenoord", 76, 1.7)
Set-Row $ws 7 @("29-01-2025 20:00", "WORLD", "UEFA CHAMPIONS LEAGUE", "Manchester City - Club Brugge KV", 83.3, 1.95)
Set-Row $ws 8 @("29-01-2025 21:30", "BRAZIL", "PAULISTA - A1", "São Bernardo - Santos", 76.7, 2.05)

# --- Sheet: Over_Under ---
$ws = $wb.Worksheets.Item("Over_Under")
Set-Row $ws 2 @("28-01-2025 19:45", "ENGLAND", "LEAGUE ONE", "Bolton - Northampton", 80, 1.7, 15, 2.75)
Set-Row $ws 3 @("28-01-2025 19:45", "ENGLAND", "LEAGUE ONE", "Exeter City - Leyton Orient", 80, 2.2, 65, 4)
Set-Row $ws 4 @("28-01-2025 19:45", "SCOTLAND", "CHAMPIONSHIP", "Hamilton Academical - Raith Rovers", 85, 1.83, 40, 3.25)
Set-Row $ws 5 @("28-01-2025 19:45", "ENGLAND", "NATIONAL LEAGUE - NORTH", "Radcliffe - Chester", 100, 1.75, 60, 2.88)
Set-Row $ws 6 @("28-01-2025 19:45", "ENGLAND", "NATIONAL LEAGUE - NORTH", "Rushall Olympic - King's Lynn Town", 80, 1.91, 5, 3.3)
Set-Row $ws 7 @("28-01-2025 15:00", "MEXICO", "U23 LEAGUE", "Puebla U23 - Mazatlan FC U23", 75, 1.8, 75, 3)
Set-Row $ws 8 @("28-01-2025 15:00", "PORTUGAL", "LIGA REVELAÇÃO U23", "Gil Vicente U23 - Vizela U23", 82.5, 1.65, 62.5, 2.6)
Set-Row $ws 9 @("29-01-2025 20:00", "WORLD", "UEFA CHAMPIONS LEAGUE", "Inter - Monaco", 70, 1.67, 65, 2.62)
Set-Row $ws 10 @("29-01-2025 20:00", "WORLD", "UEFA CHAMPIONS LEAGUE", "VfB Stuttgart - Paris Saint Germain", 80, 2, 53.3, 3.2)
